# Transportation file updates from CB
$wb = $excel.ActiveWorkbook

# --- AVLo-freight: B7 now tracks the vehicle-weight input in B3 instead of
# being hard-coded to 0. The shared formulas in C7:AJ7 already reference
# $B7, so they recompute automatically once B7 becomes a formula.
$wsFreight = $wb.Worksheets.Item("AVLo-freight")
$wsFreight.Range("B7").Formula = "=B3"

# --- Make AVLo-freight the active/selected sheet (it was "AVLo-passengers"
# before), with B8 selected.
$wsFreight.Activate() | Out-Null
$wsFreight.Range("B8").Select() | Out-Null
